# Mise à jour de l'application
# Appends a new training-session ("Entrainement") block of 9 player rows
# (one MD J+3 session dated 2025-10-07, serial 45937) to the bottom of the
# existing GPS data log on Feuil1, mirroring the layout of the rows already
# present (same Type/Période/MD columns, one row per player).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Last existing data row is 627; new rows go to 628..636.
$firstNewRow = 628
$lastExistingRow = 627
$lastNewRow = 636

# Clone the number formats / styles of the last existing row (B = date
# format, D = centered text) onto the new block, exactly like copy/pasting
# the row down and then overwriting the values.
$ws.Range("A" + $lastExistingRow + ":V" + $lastExistingRow).Copy() | Out-Null
$ws.Range("A" + $firstNewRow + ":V" + $lastNewRow).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$rows = @(
    @{ Row=628; E="Yoann Martelat";   F="center midfield"; G="01:40:03"; H=7.58; I=0.87; J=6.7;  K=0.84; L=0.04; M=0;    N=0; O=0;  P=4.17; Q=22.08; R=3.7;  S=32; T=0;  U=17; V=3  },
    @{ Row=629; E="Ilan Ihaddadene";  F="center midfield"; G="01:40:13"; H=8.61; I=1.08; J=7.53; K=1.06; L=0.03; M=0;    N=0; O=0;  P=5.07; Q=20.94; R=4.39; S=71; T=7;  U=39; V=4  },
    @{ Row=630; E="Omar Benyounes";   F="center midfield"; G="01:39:55"; H=8.71; I=1.12; J=7.58; K=1.07; L=0.06; M=0;    N=0; O=0;  P=5.16; Q=22.72; R=4.49; S=61; T=8;  U=41; V=1  },
    @{ Row=631; E="Ilyes Boughanmi";  F="center forward";  G="01:39:56"; H=7.94; I=0.83; J=7.09; K=0.77; L=0.08; M=0;    N=0; O=0;  P=3.97; Q=24.16; R=5.08; S=39; T=6;  U=42; V=11 },
    @{ Row=632; E="Emmanuel Valey";   F="left forward";    G="01:37:32"; H=7.65; I=0.26; J=7.38; K=0.27; L=0.01; M=0;    N=0; O=0;  P=3.99; Q=22.28; R=5.48; S=64; T=8;  U=46; V=9  },
    @{ Row=633; E="Karahali Souaré";  F="right forward";   G="01:40:39"; H=7.55; I=0.46; J=7.08; K=0.39; L=0.06; M=0.03; N=0; O=3;  P=3.84; Q=30.03; R=6.27; S=43; T=21; U=39; V=21 },
    @{ Row=634; E="Malik Boussaid";   F="right back";      G="01:39:46"; H=8.71; I=1.21; J=7.49; K=0.74; L=0.49; M=0;    N=0; O=2;  P=4.57; Q=25.71; R=5.41; S=49; T=17; U=40; V=18 },
    @{ Row=635; E="Mattheo Haon";     F="right back";      G="01:38:08"; H=8.63; I=1.15; J=7.47; K=1.1;  L=0.07; M=0;    N=0; O=0;  P=5.21; Q=24.09; R=4.7;  S=54; T=7;  U=21; V=5  },
    @{ Row=636; E="Kamal Bafounta";   F="center midfield"; G="01:40:13"; H=7.67; I=1.02; J=6.63; K=0.93; L=0.1;  M=0;    N=0; O=0;  P=4.48; Q=24.17; R=4.3;  S=54; T=3;  U=39; V=7  }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = "Entrainement"     # A - Type
    $ws.Cells.Item($row, 2).Value = 45937              # B - Date (2025-10-07)
    $ws.Cells.Item($row, 3).Value = "Global"            # C - Période
    $ws.Cells.Item($row, 4).Value = "J+3"               # D - MD
    $ws.Cells.Item($row, 5).Value = $r.E                # E - Nom du joueur
    $ws.Cells.Item($row, 6).Value = $r.F                # F - Poste
    $ws.Cells.Item($row, 7).Value = $r.G                # G - Temps joué
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
    $ws.Cells.Item($row, 21).Value = $r.U
    $ws.Cells.Item($row, 22).Value = $r.V
}

# Match the refreshed view state from the diff: selection sits just past
# the newly appended block, and the visible window has scrolled down with
# it (best-effort — host may not persist scroll position).
$ws.Range("C639").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 603
